# Journal de travail James - add new journal entries (rows 14-17) and
# update the Total row's formula result accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily clear the "Total" label so that, when re-entered after the
# new activity strings below, it lands at the end of the shared-string table.
$ws.Range("B32").Value = ""

# New journal entries (Date, Activité, Heures)
$ws.Range("A14").Value = (Get-Date -Year 2018 -Month 3 -Day 26 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("B14").Value = "Discussion, planification et organisation de groupe"
$ws.Range("C14").Value = 1.5

$ws.Range("A15").Value = (Get-Date -Year 2018 -Month 3 -Day 26 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("B15").Value = "Conception et analyse de la drawZone"
$ws.Range("C15").Value = 3

$ws.Range("A16").Value = (Get-Date -Year 2018 -Month 4 -Day 9 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("B16").Value = "Brainstrom et résolution de bug graphique"
$ws.Range("C16").Value = 1.5

$ws.Range("A17").Value = (Get-Date -Year 2018 -Month 4 -Day 15 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("B17").Value = "Création du projet sur IceScrum (prend trop de temps : abondon)"
$ws.Range("C17").Value = 1

# Row heights reflow to fit the new content (rows 14-16 are single line,
# row 17 wraps to a taller line because of its longer activity text)
$ws.Rows.Item(14).RowHeight = 13.8
$ws.Rows.Item(15).RowHeight = 13.8
$ws.Rows.Item(16).RowHeight = 13.8
$ws.Rows.Item(17).RowHeight = 23.85

# Restore the "Total" label (re-added after the new strings above)
$ws.Range("B32").Value = "Total"

# Move the active selection to G7 (as recorded in the saved view state)
$ws.Range("G7").Select()
